$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regimen rows appended below the existing table (rows 30-34),
# following the same id/displayName column layout as the rest of the sheet.
$ws.Range("A30").Value = "PTV_r_OBV_DSV"
$ws.Range("B30").Value = "PTV/r/OBV/DSV"

$ws.Range("A31").Value = "PTV_r_OBV_DSV_RBV"
$ws.Range("B31").Value = "PTV/r/OBV/DSV/RBV"

$ws.Range("A32").Value = "PTV_r_OBV_RBV"
$ws.Range("B32").Value = "PTV/r/OBV/RBV"

$ws.Range("A33").Value = "PTV_r_OBV"
$ws.Range("B33").Value = "PTV/r/OBV"

$ws.Range("A34").Value = "PIB_PTV_r_RBV"
$ws.Range("B34").Value = "PIB/PTV/r/RBV"

# A31:A33 pick up the same font/alignment formatting used for the last
# row of the previous block (A29), but without any border.
$ws.Range("A29").Copy()
$ws.Range("A31:A33").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A31:A33").Borders.LineStyle = -4142  # xlLineStyleNone
$excel.CutCopyMode = $false

# Scroll the view so the newly-added rows are visible, mirroring the
# author's final cursor position in the saved workbook.
$ws.Range("A24").Select()
